$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "stars" column (E) entries that were missing comma separators between names
$ws.Range("E17").Value = 'Amy Poehler, Bill Hader, Lewis Black'
$ws.Range("E15").Value = 'Jesse Eisenberg, Andrew Garfield, Justin Timberlake'

# Normalize the "year_genre" column (H) separator from the bullet character "•" to a pipe "|"
$ws.Range("H2").Value  = '2019 | Comedy/Western | 2h 40m'
$ws.Range("H3").Value  = '2022 | Action/Crime | 2h 56m'
$ws.Range("H4").Value  = '2010 | Horror/Drama | 1h 48m'
$ws.Range("H5").Value  = '2010 | Sci-fi/Action | 2h 28m'
$ws.Range("H6").Value  = '2023 | Comedy/Fantasy | 1h 54m'
$ws.Range("H7").Value  = '2023 | Drama/History | 3h'
$ws.Range("H8").Value  = '2019 | Comedy/Romance | 2h 5m'
$ws.Range("H9").Value  = '2018 | Sci-fi | 2 seasons'
$ws.Range("H10").Value = '2006 | Thriller/Sci-fi | 2h 10m'
$ws.Range("H11").Value = '2014 | Sci-fi/Adventure | 2h 49m'
$ws.Range("H12").Value = '2008 | Comedy/Romance | 1h 51m'
$ws.Range("H13").Value = '2009 | Comedy/Romance | 1h 45m'
$ws.Range("H14").Value = '2012 | Romance/Comedy | 1h 32m'
$ws.Range("H15").Value = '2010 | Drama/Historical drama | 2 hours'
$ws.Range("H16").Value = '2019 | Family/Comedy | 1h 42m'
$ws.Range("H17").Value = '2015 | Family/Comedy | 1h 35m'

# Move the active selection to H8, matching the workbook's saved cursor position
$ws.Range("H8").Select()
